# Refresh cryptos list data (prices / 1h volume, and row realignment)
# as published by the upstream GitHub Actions scraper job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '27.126.17'
$ws.Range("E2").Value = '  -0.59%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '1.895.31'
$ws.Range("E3").Value = '  -0.69%  '
# Row 4: TetherUSD
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.07%  '
# Row 5: BNB
$ws.Range("D5").Value = '''307.11'
$ws.Range("E5").Value = '  -0.25%  '
# Row 6: USDC
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.13%  '
# Row 7: XRP
$ws.Range("D7").Value = '''0.5209'
$ws.Range("E7").Value = '  -0.80%  '
# Row 8: Cardano
$ws.Range("D8").Value = '''0.3764'
$ws.Range("E8").Value = '  -0.68%  '
# Row 9: Dogecoin
$ws.Range("D9").Value = '''0.07273'
$ws.Range("E9").Value = '  +0.09%  '
# Row 10: Solana
$ws.Range("D10").Value = '''21.16'
$ws.Range("E10").Value = '  -0.86%  '
# Row 11: Polygon
$ws.Range("D11").Value = '''0.9007'
$ws.Range("E11").Value = '  -0.16%  '
# Row 12: TRON
$ws.Range("D12").Value = '''0.08215'
$ws.Range("E12").Value = '  -0.88%  '
# Row 13: WrappedEther
$ws.Range("D13").Value = '1.978.74'
$ws.Range("E13").Value = '  +3.65%  '
# Row 14: Litecoin
$ws.Range("D14").Value = '''96.47'
$ws.Range("E14").Value = '  +1.24%  '
# Row 15: Polkadot
$ws.Range("D15").Value = '''5.308'
$ws.Range("E15").Value = '  +0.23%  '
# Row 17: ShibaInu
$ws.Range("D17").Value = '''0.000008613'
$ws.Range("E17").Value = '  -0.08%  '
# Row 18: Avalanche
$ws.Range("D18").Value = '''14.61'
$ws.Range("E18").Value = '  +0.77%  '
# Row 19: Dai
$ws.Range("E19").Value = '  +0.26%  '
# Row 20: WrappedBTC
$ws.Range("D20").Value = '27.156.59'
$ws.Range("E20").Value = '  -0.65%  '
# Row 21: Uniswap
$ws.Range("D21").Value = '''5.082'
$ws.Range("E21").Value = '  +0.10%  '
# Row 22: Cosmos
$ws.Range("B22").Value = 'Cosmos'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D22").Value = '''10.71'
$ws.Range("E22").Value = '  +0.46%  '
# Row 23: Chainlink
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '''6.424'
$ws.Range("E23").Value = '  -0.71%  '
# Row 24: Monero
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '''148.48'
$ws.Range("E24").Value = '  +1.49%  '
# Row 25: LidoDAOToken
$ws.Range("D25").Value = '''2.311'
$ws.Range("E25").Value = '  +0.03%  '
# Row 26: EthereumClassic
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '''18.21'
$ws.Range("E26").Value = '  +0.08%  '
# Row 27: Toncoin
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '''1.742'
$ws.Range("E27").Value = '  -0.35%  '
# Row 28: BitcoinCash
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '''115.29'
$ws.Range("E28").Value = '  +0.20%  '
# Row 29: InternetComputer(DFINITY)
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''4.804'
$ws.Range("E29").Value = '  -0.22%  '
# Row 30: Filecoin
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''4.863'
$ws.Range("E30").Value = '  -2.70%  '
# Row 31: Stellar
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.09205'
$ws.Range("E31").Value = '  -0.32%  '
# Row 32: Hedera
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.05020'
$ws.Range("E32").Value = '  -0.81%  '
# Row 33: ImmutableX
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7927'
$ws.Range("E33").Value = '  -1.79%  '
# Row 34: ARBITRUM
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.219'
$ws.Range("E34").Value = '  -2.04%  '
# Row 35: MXToken
$ws.Range("B35").Value = 'MXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D35").Value = '''3.430'
$ws.Range("E35").Value = '  +1.66%  '
# Row 36: HuobiToken
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '''2.969'
$ws.Range("E36").Value = '  -0.27%  '
# Row 37: RenderToken
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '''2.609'
$ws.Range("E37").Value = '  +1.22%  '
# Row 38: TheSandbox
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '''0.5725'
$ws.Range("E38").Value = '  -0.12%  '
# Row 39: VeChain
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.02002'
$ws.Range("E39").Value = '  +0.83%  '
# Row 40: TrustWalletToken
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''1.075'
$ws.Range("E40").Value = '  -0.24%  '
# Row 41: Aptos
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '''9.028'
$ws.Range("E41").Value = '  +0.06%  '
# Row 42: FraxShare
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''6.569'
$ws.Range("E42").Value = '  -0.90%  '
# Row 43: Quant
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''116.50'
$ws.Range("E43").Value = '  -2.66%  '
# Row 44: Algorand
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '''0.1515'
$ws.Range("E44").Value = '  -0.24%  '
# Row 45: Decentraland
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.4875'
$ws.Range("E45").Value = '  +0.66%  '
# Row 46: PaxDollar
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''1.002'
$ws.Range("E46").Value = '  +0.15%  '
# Row 47: EnergySwap
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''10.09'
$ws.Range("E47").Value = '  -1.07%  '
# Row 48: NEARProtocol
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.623'
$ws.Range("E48").Value = '  +0.24%  '
# Row 49: Elrond
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''38.30'
$ws.Range("E49").Value = '  +1.60%  '
# Row 50: Aave
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '''63.80'
$ws.Range("E50").Value = '  -0.09%  '
# Row 51: Cronos
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.05927'
$ws.Range("E51").Value = '  -0.49%  '
